# Applies updated Betfair Back/Lay odds for 2026-01-06 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.85
$ws.Range("G2").Value = 1.86
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 3.85
$ws.Range("N2").Value = 3.4
$ws.Range("W2").Value = 2.16
$ws.Range("AB2").Value = 7.8
# Row 3
$ws.Range("I3").Value = 1.83
$ws.Range("J3").Value = 3.65
$ws.Range("P3").Value = 1.79
$ws.Range("U3").Value = 1.88
$ws.Range("V3").Value = 2.2
$ws.Range("Y3").Value = 7.6
$ws.Range("Z3").Value = 9.800000000000001
$ws.Range("AC3").Value = 8
$ws.Range("AM3").Value = 150
$ws.Range("AO3").Value = 14
# Row 4
$ws.Range("H4").Value = 1.79
$ws.Range("I4").Value = 1.8
$ws.Range("V4").Value = 2.24
$ws.Range("X4").Value = 9.800000000000001
$ws.Range("AM4").Value = 200
$ws.Range("AO4").Value = 16.5
# Row 5
$ws.Range("F5").Value = 5.4
$ws.Range("G5").Value = 5.5
$ws.Range("L5").Value = 1.42
$ws.Range("Q5").Value = 2.06
# Row 6
$ws.Range("F6").Value = 2.76
$ws.Range("G6").Value = 2.84
$ws.Range("H6").Value = 2.72
$ws.Range("I6").Value = 2.78
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 3.95
$ws.Range("U6").Value = 2.12
$ws.Range("V6").Value = 1.56
$ws.Range("W6").Value = 1.54
$ws.Range("AA6").Value = 50
$ws.Range("AD6").Value = 12.5
# Row 7
$ws.Range("F7").Value = 3.35
$ws.Range("G7").Value = 3.4
$ws.Range("H7").Value = 2.36
$ws.Range("I7").Value = 2.4
$ws.Range("L7").Value = 1.42
$ws.Range("P7").Value = 1.95
$ws.Range("Q7").Value = 2.02
$ws.Range("S7").Value = 3.6
$ws.Range("T7").Value = 1.79
$ws.Range("U7").Value = 2.2
$ws.Range("V7").Value = 1.71
$ws.Range("W7").Value = 1.41
$ws.Range("Z7").Value = 14.5
$ws.Range("AB7").Value = 13.5
$ws.Range("AI7").Value = 38
$ws.Range("AN7").Value = 36
# Row 8
$ws.Range("H8").Value = 6.4
$ws.Range("K8").Value = 4.3
$ws.Range("P8").Value = 2.28
$ws.Range("Q8").Value = 1.71
$ws.Range("S8").Value = 2.68
$ws.Range("T8").Value = 1.77
$ws.Range("U8").Value = 2.16
$ws.Range("V8").Value = 1.16
$ws.Range("X8").Value = 20
$ws.Range("AE8").Value = 85
$ws.Range("AF8").Value = 9.800000000000001
$ws.Range("AH8").Value = 19
$ws.Range("AK8").Value = 15
